$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $text)
    $range.NumberFormat = "@"
    if ($text.Length -eq 0) {
        $range.Value = ""
    } else {
        $range.Value = ("X" * $text.Length)
        for ($i = 0; $i -lt $text.Length; $i++) {
            $range.Characters($i+1, 1).Text = $text.Substring($i,1)
        }
    }
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "35.405.34"
$ws.Range("E2").Value = "  +2.68%  "

Set-TextValue $ws.Range("D3") "1.841.02"
$ws.Range("E3").Value = "  +1.83%  "

$ws.Range("E4").Value = "  +0.44%  "

Set-TextValue $ws.Range("D5") "230.72"
$ws.Range("E5").Value = "  +2.72%  "

Set-TextValue $ws.Range("D6") "0.608"
$ws.Range("E6").Value = "  +1.26%  "

$ws.Range("E7").Value = "  +0.49%  "

Set-TextValue $ws.Range("D8") "43.23"
$ws.Range("E8").Value = "  +11.18%  "

Set-TextValue $ws.Range("D9") "0.308"
$ws.Range("E9").Value = "  +7.16%  "

Set-TextValue $ws.Range("D10") "0.0700"
$ws.Range("E10").Value = "  +4.41%  "

$ws.Range("E11").Value = "  +3.74%  "

Set-TextValue $ws.Range("D12") "2.107.22"
$ws.Range("E12").Value = "  +1.89%  "

Set-TextValue $ws.Range("D13") "1.841.38"
$ws.Range("E13").Value = "  +1.93%  "

Set-TextValue $ws.Range("D14") "11.27"
$ws.Range("E14").Value = "  +1.23%  "

Set-TextValue $ws.Range("D15") "0.671"
$ws.Range("E15").Value = "  +6.69%  "

$ws.Range("E16").Value = "  +6.56%  "

Set-TextValue $ws.Range("D17") "35.380.58"
$ws.Range("E17").Value = "  +2.71%  "

Set-TextValue $ws.Range("D18") "69.98"
$ws.Range("E18").Value = "  +2.63%  "

Set-TextValue $ws.Range("D19") [string]::Concat("0.0", [char]0x2083, "0798")
$ws.Range("E19").Value = "  +3.88%  "

Set-TextValue $ws.Range("D20") "244.60"
$ws.Range("E20").Value = "  +1.30%  "

$ws.Range("E21").Value = "  +8.55%  "

Set-TextValue $ws.Range("D22") "4.68"
$ws.Range("E22").Value = "  +14.20%  "

$ws.Range("E23").Value = "  +0.31%  "

$ws.Range("E24").Value = "  +0.75%  "

Set-TextValue $ws.Range("D25") "169.24"
$ws.Range("E25").Value = "  -0.99%  "

$ws.Range("E26").Value = "  +2.59%  "

Set-TextValue $ws.Range("D27") "17.71"
$ws.Range("E27").Value = "  +1.05%  "

$ws.Range("E28").Value = "  +0.59%  "

Set-TextValue $ws.Range("D29") "1.55"
$ws.Range("E29").Value = "  +26.60%  "

Set-TextValue $ws.Range("D30") "1.01"
$ws.Range("E30").Value = "  +0.46%  "

Set-TextValue $ws.Range("D31") "3.266.64"
$ws.Range("E31").Value = "  +34.45%  "

$ws.Range("E32").Value = "  +6.20%  "

$ws.Range("E33").Value = "  +4.29%  "

Set-TextValue $ws.Range("D34") "4.05"
$ws.Range("E34").Value = "  +5.78%  "

$ws.Range("E35").Value = "  +1.15%  "

Set-TextValue $ws.Range("D36") "94.05"
$ws.Range("E36").Value = "  +14.06%  "

Set-TextValue $ws.Range("D37") "0.685"
$ws.Range("E37").Value = "  +6.96%  "

Set-TextValue $ws.Range("D38") "1.342.71"
$ws.Range("E38").Value = "  +2.07%  "

$ws.Range("E39").Value = "  +2.65%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D40") "2.43"
$ws.Range("E40").Value = "  +5.42%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D41") "0.0194"
$ws.Range("E41").Value = "  +4.48%  "

Set-TextValue $ws.Range("D42") "15.19"
$ws.Range("E42").Value = "  +10.56%  "

Set-TextValue $ws.Range("D43") "1.00"
$ws.Range("E43").Value = "  +6.18%  "

$ws.Range("E44").Value = "  +4.42%  "

$ws.Range("E45").Value = "  +1.01%  "

$ws.Range("E46").Value = "  +0.34%  "

Set-TextValue $ws.Range("D47") "6.23"
$ws.Range("E47").Value = "  +7.57%  "

Set-TextValue $ws.Range("D48") "0.0519"
$ws.Range("E48").Value = "  +1.29%  "

Set-TextValue $ws.Range("D49") "2.008.76"
$ws.Range("E49").Value = "  +2.04%  "

$ws.Range("E50").Value = "  +0.55%  "

